$d = $word.ActiveDocument

# The "BODY PART – " Mad Lib blank is where the cursor was last left in
# the document (marked by the built-in "_GoBack" bookmark). Put the
# insertion point there and type the answer, same as a user clicking
# into that spot and filling in the blank.
if ($d.Bookmarks.Exists("_GoBack")) {
    $ip = $d.Bookmarks.Item("_GoBack").Range
} else {
    $ip = $d.Content
    $ip.Find.Execute("BODY PART – ", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0)
}
$ip.Collapse(0)
$ip.InsertAfter("Face muscles")

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
